# Refresh the cryptocurrency price/volume snapshot cells (columns D and E)
# with the latest scraped values, keeping them as plain text so formats
# like "94.00" / "0.000008790" / "27.300.35" survive intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.300.35'
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = '1.832.48'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("D4").Value = '''1.012'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.99%  '
$ws.Range("D5").Value = '''314.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.60%  '
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("D7").Value = '''0.4735'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.72%  '
$ws.Range("D8").Value = '''0.3687'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.60%  '
$ws.Range("D9").Value = '''0.07438'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("D10").Value = '''0.8861'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.60%  '
$ws.Range("D11").Value = '''20.49'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("D12").Value = '1.886.36'
$ws.Range("E12").Value = '  +4.53%  '
$ws.Range("D13").Value = '''0.07337'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.18%  '
$ws.Range("D14").Value = '''5.427'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").Value = '''94.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.94%  '
$ws.Range("D16").Value = '''6.558'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").Value = '''0.000008790'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.21%  '
$ws.Range("D19").Value = '''1.011'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = '27.514.10'
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("E21").Value = '  +0.88%  '
$ws.Range("E23").Value = '  +0.72%  '
$ws.Range("D24").Value = '2.097.14'
$ws.Range("E24").Value = '  +2.97%  '
$ws.Range("D25").Value = '''1.894'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = '''152.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("E27").Value = '  +1.14%  '
$ws.Range("D28").Value = '''2.143'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("D29").Value = '''5.229'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("D30").Value = '''116.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("D31").Value = '''0.08989'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.04%  '
$ws.Range("D32").Value = '''0.7497'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.19%  '
$ws.Range("D33").Value = '''1.174'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.88%  '
$ws.Range("E34").Value = '  +1.33%  '
$ws.Range("D35").Value = '''2.941'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.37%  '
$ws.Range("E36").Value = '  +0.90%  '
$ws.Range("D37").Value = '''1.094'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").Value = '''0.05346'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.06%  '
$ws.Range("D39").Value = '''0.01951'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.28%  '
$ws.Range("D40").Value = '''2.976'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").Value = '''2.403'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.31%  '
$ws.Range("D42").Value = '''7.234'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.16%  '
$ws.Range("D43").Value = '''0.5288'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").Value = '''8.475'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("D46").Value = '''0.4928'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.63%  '
$ws.Range("D47").Value = '''10.55'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.29%  '
$ws.Range("E48").Value = '  +0.84%  '
$ws.Range("D49").Value = '''104.86'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.52%  '
$ws.Range("E50").Value = '  +0.54%  '
$ws.Range("D51").Value = '''0.06302'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.21%  '
